$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by exactly one day
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 1
}

# Update the Actual Production (MW) values for rows 24-41 (column B) to reflect
# the new day's solar production figures
$newB = @{
    24 = 4
    25 = 23
    26 = 56
    27 = 105
    28 = 182
    29 = 280
    30 = 390
    31 = 487
    32 = 584
    33 = 700
    34 = 816
    35 = 916
    36 = 998
    37 = 1036
    38 = 1110
    39 = 1173
    40 = 1233
    41 = 0
}

foreach ($row in $newB.Keys) {
    $ws.Cells.Item($row, 2).Value = $newB[$row]
}
